$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.833.19"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.695.96"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.46%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.02"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.72"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.607"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.728.21"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.69"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.113"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +6.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.387"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.157"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.177.98"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.46"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +8.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.733.62"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000151"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +7.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.713.06"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.90"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.86"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "362.30"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.03"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.58%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.533"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.45"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.166"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.58"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +7.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.993"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.04"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +7.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0848"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +5.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.13"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +11.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "169.94"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.997"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.48"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.17"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +18.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.75"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +11.00%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +6.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.82"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +11.04%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +20.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "350.54"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +12.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.23"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +8.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.94"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.70"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +14.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.54"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +8.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0592"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +7.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.72"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +8.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0258"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +6.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.638"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "137.00"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.42%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.24%  "
